$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
$cs = $nm.ColorScheme
$cs.Colors(1).RGB = 255
Write-Output "done"
